# Generate Report for Handback
# Update the Handoff/Handback timestamp cells for the
# "ee81f607-5d83-48c5-baf3-4679802b6fad" file across the Overview,
# zh-cn and de-de worksheets, as produced by a new report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for ee81f607-5d83-48c5-baf3-4679802b6fad.md
$wsOverview.Range("G4").Value = "2016-08-18 20:47:05"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-08-18 20:46:55"
$wsZhCn.Range("K4").Value = "2016-08-18 20:47:34"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H4").Value = "2016-08-18 20:47:05"
$wsDeDe.Range("K4").Value = "2016-08-18 20:47:42"
